# Update the mentors addendum partnership-term sentence from
# "5 years of partnership (Fall 2025 through Spring 2030)" to
# "5 years of partnership (Fall 2026 through Spring 2031)".

$d = $word.ActiveDocument

# Bump the start year of the partnership term: Fall 2025 -> Fall 2026
$d.Content.Find.Execute(
    "Fall 2025", $true, $false, $false, $false, $false, $true, 1, $false,
    "Fall 2026", 2
)

# Bump the end year of the partnership term: Spring 2030 -> Spring 2031
$d.Content.Find.Execute(
    "Spring 2030", $true, $false, $false, $false, $false, $true, 1, $false,
    "Spring 2031", 2
)
